# Add season record columns (Wins / Losses / Ties) to the roster sheet.
# The previous export only captured team statistics, not the season
# win/loss/tie record, so we extend the table with three new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44

# --- Header row (row 1) -----------------------------------------------
# Copy the formatting of the last existing header cell (AC1) onto the
# three new header cells so they keep the bold / bordered / centered
# header style, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-44) ---------------------------------------------
# Every player row shares the same season record for this team/year.
$wins = 87
$losses = 75
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
